$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.718.79'
$ws.Range('E2').Value = '  -1.68%  '
$ws.Range('D3').Value = '2.074.75'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.55'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.58'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.395'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').Value = '2.379.36'
$ws.Range('E12').Value = '  -2.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.79'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.94'
$ws.Range('E14').Value = '  -3.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.773'
$ws.Range('E15').Value = '  -2.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.37'
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('D17').Value = '2.094.96'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').Value = '37.638.68'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.43'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '0.0₃0834'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.17'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.36'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.47'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('E31').Value = '  +1.61%  '
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  -5.69%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.35'
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0975'
$ws.Range('E40').Value = '  -3.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.58'
$ws.Range('E41').Value = '  +2.00%  '
$ws.Range('E42').Value = '  -2.55%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.71'
$ws.Range('E44').Value = '  +5.64%  '
$ws.Range('D45').Value = '1.436.24'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('E46').Value = '  -1.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.20'
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.39'
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').Value = '2.264.40'
$ws.Range('E51').Value = '  -2.31%  '
